$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New registrant rows, mirroring the formatting quirks produced by the
# CSV-driven registration import (row 2 was the first imported record;
# rows 3-6 are the newly imported records).
$rows = @(
    @{ Row = 3; C = 5;  D = 11; Email = "sonytuladhar25@gmail.com"; Pwd = "Khalti1" },
    @{ Row = 4; C = 6;  D = 12; Email = "sonytuladhar26@gmail.com"; Pwd = "Khalti2" },
    @{ Row = 5; C = 7;  D = 13; Email = "sonytuladhar27@gmail.com"; Pwd = "Khalti3" },
    @{ Row = 6; C = 8;  D = 14; Email = "sonytuladhar28@gmail.com"; Pwd = "Khalti4" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy the whole style/formatting of row 2 down onto the new row first,
    # so every cell lands on the same base style (s="1") that row 2 uses.
    $ws.Range("A2:I2").Copy()
    $ws.Range("A" + $rowNum + ":I" + $rowNum).PasteSpecial(-4122)

    # C, D, F and G use the alternate style (s="2") seen on E1 in the
    # original sheet, matching the target layout.
    $ws.Range("E1").Copy()
    $ws.Range("C" + $rowNum).PasteSpecial(-4122)
    $ws.Range("D" + $rowNum).PasteSpecial(-4122)
    $ws.Range("F" + $rowNum).PasteSpecial(-4122)
    $ws.Range("G" + $rowNum).PasteSpecial(-4122)

    $ws.Range("A" + $rowNum).Value = "Sony"
    $ws.Range("B" + $rowNum).Value = "Tuladhar"
    $ws.Range("C" + $rowNum).Value = $r.C
    $ws.Range("D" + $rowNum).Value = $r.D
    $ws.Range("E" + $rowNum).Value = 1994
    $ws.Range("F" + $rowNum).Value = $r.Email
    $ws.Range("G" + $rowNum).Value = $r.Pwd
    $ws.Range("H" + $rowNum).Value = 76697669
    $ws.Range("I" + $rowNum).Value = 76697669
}

$excel.CutCopyMode = 0
